$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as TEXT, preserving default (unstyled) formatting like the source file.
function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextCell 2 4 '68.119.00'
$ws.Cells.Item(2, 5).Value = '  -0.37%  '
Set-TextCell 3 4 '3.662.99'
$ws.Cells.Item(3, 5).Value = '  -1.00%  '
Set-TextCell 4 4 '0.998'
$ws.Cells.Item(4, 5).Value = '  -0.24%  '
Set-TextCell 5 4 '595.91'
$ws.Cells.Item(5, 5).Value = '  +1.87%  '
Set-TextCell 6 4 '192.60'
$ws.Cells.Item(6, 5).Value = '  +5.17%  '
Set-TextCell 7 4 '0.621'
$ws.Cells.Item(7, 5).Value = '  -1.12%  '
$ws.Cells.Item(8, 5).Value = '  -0.09%  '
Set-TextCell 9 4 '0.700'
$ws.Cells.Item(9, 5).Value = '  -2.51%  '
$ws.Cells.Item(10, 5).Value = '  -6.97%  '
Set-TextCell 11 4 '56.84'
$ws.Cells.Item(11, 5).Value = '  +2.16%  '
Set-TextCell 12 4 '0.0000271'
$ws.Cells.Item(12, 5).Value = '  -7.07%  '
Set-TextCell 13 4 '10.19'
$ws.Cells.Item(13, 5).Value = '  -2.01%  '
Set-TextCell 14 4 '4.236.66'
$ws.Cells.Item(14, 5).Value = '  +1.09%  '
Set-TextCell 15 4 '3.650.51'
$ws.Cells.Item(15, 5).Value = '  -1.39%  '
$ws.Cells.Item(16, 5).Value = '  +0.38%  '
Set-TextCell 17 4 '18.86'
$ws.Cells.Item(17, 5).Value = '  -2.83%  '
$ws.Cells.Item(18, 2).Value = 'Polygon'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell 18 4 '1.11'
$ws.Cells.Item(18, 5).Value = '  -1.68%  '
$ws.Cells.Item(19, 2).Value = 'WrappedBTC'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell 19 4 '67.813.33'
$ws.Cells.Item(19, 5).Value = '  -0.45%  '
Set-TextCell 20 4 '12.42'
$ws.Cells.Item(20, 5).Value = '  -2.67%  '
Set-TextCell 21 4 '399.86'
$ws.Cells.Item(21, 5).Value = '  -2.27%  '
Set-TextCell 22 4 '4.37'
$ws.Cells.Item(22, 5).Value = '  -2.40%  '
Set-TextCell 23 4 '87.76'
$ws.Cells.Item(23, 5).Value = '  -0.90%  '
Set-TextCell 24 4 '2.95'
$ws.Cells.Item(24, 5).Value = '  -2.45%  '
Set-TextCell 25 4 '11.04'
$ws.Cells.Item(25, 5).Value = '  -0.31%  '
Set-TextCell 26 4 '12.49'
$ws.Cells.Item(26, 5).Value = '  -2.35%  '
Set-TextCell 27 4 '6.07'
$ws.Cells.Item(27, 5).Value = '  +0.21%  '
Set-TextCell 28 4 '3.65'
$ws.Cells.Item(28, 5).Value = '  -6.41%  '
Set-TextCell 29 4 '9.33'
$ws.Cells.Item(29, 5).Value = '  -1.70%  '
Set-TextCell 30 4 '31.81'
$ws.Cells.Item(30, 5).Value = '  -2.68%  '
Set-TextCell 31 4 '7.29'
$ws.Cells.Item(31, 5).Value = '  -1.66%  '
Set-TextCell 32 4 '12.28'
$ws.Cells.Item(32, 5).Value = '  -1.93%  '
Set-TextCell 33 4 '44.31'
$ws.Cells.Item(33, 5).Value = '  +1.53%  '
Set-TextCell 34 4 '65.78'
$ws.Cells.Item(34, 5).Value = '  +0.47%  '
$ws.Cells.Item(35, 2).Value = 'Hedera'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 35 4 '0.115'
$ws.Cells.Item(35, 5).Value = '  -1.49%  '
$ws.Cells.Item(36, 2).Value = 'Bittensor'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell 36 4 '603.97'
$ws.Cells.Item(36, 5).Value = '  +1.48%  '
$ws.Cells.Item(37, 5).Value = '  +0.05%  '
Set-TextCell 38 4 '0.393'
$ws.Cells.Item(38, 5).Value = '  -2.26%  '
$ws.Cells.Item(39, 5).Value = '  -0.18%  '
Set-TextCell 40 4 '0.0₃0767'
$ws.Cells.Item(40, 5).Value = '  -13.88%  '
$ws.Cells.Item(41, 5).Value = '  -0.93%  '
Set-TextCell 42 4 '2.87'
$ws.Cells.Item(42, 5).Value = '  -4.12%  '
Set-TextCell 43 4 '0.0425'
$ws.Cells.Item(43, 5).Value = '  -2.35%  '
Set-TextCell 44 4 '2.53'
$ws.Cells.Item(44, 5).Value = '  -8.40%  '
$ws.Cells.Item(45, 5).Value = '  +1.04%  '
Set-TextCell 46 4 '2.769.34'
$ws.Cells.Item(46, 5).Value = '  -0.17%  '
$ws.Cells.Item(47, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell 47 4 '3.12'
$ws.Cells.Item(47, 5).Value = '  -0.34%  '
$ws.Cells.Item(48, 2).Value = 'Monero'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 48 4 '143.39'
$ws.Cells.Item(48, 5).Value = '  +2.88%  '
Set-TextCell 49 4 '8.71'
$ws.Cells.Item(49, 5).Value = '  -6.31%  '
$ws.Cells.Item(50, 5).Value = '  -3.50%  '
Set-TextCell 51 4 '2.49'
$ws.Cells.Item(51, 5).Value = '  -15.26%  '
